$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: correct the last-checked timestamp for rows 870-883 -------------
# (tiny re-check correction baked into this "actualizar" run, before the
# 14 new rows for the 02-06-2021 00:03:54 check are appended)
for ($r = 870; $r -le 883; $r++) {
    $ws.Cells.Item($r, 4).Value = 44232.98146739583
}

# --- Step 2: append the new batch of 14 monitored-service rows --------------
# Same 14-row cycle (Nombre / URL / Disponibilidad / Fecha) used throughout
# the sheet, now stamped with the new check timestamp.
$names = @("Odoo", "Blackbox", "PowerBI", "Dropbox", "Odoo", "GEE", "UtilidadesOdoo", "Filtros Dashboard", "MapStore", "GeoServer", "Tomcat", "Shiny", "Github", "EZ Exporter")
$urls  = @(
    "https://www.dataintelligence-group.com/",
    "https://serviciodashboard.azurewebsites.net/",
    "https://powerbi.microsoft.com/es-es/",
    "https://www.dropbox.com/",
    "https://dataintelligence.store/",
    "https://app-data-i.users.earthengine.app/",
    "https://odooutil.azurewebsites.net/",
    "https://filtradordashboard.azurewebsites.net/",
    "https://ide.dataintelligence-group.com/mapstore/#/",
    "https://ide.dataintelligence-group.com/geoserver/web/?0",
    "https://ide.dataintelligence-group.com/",
    "https://rpubs.com/dataintelligence/",
    "https://github.com/Sud-Austral/",
    "https://ezexporter.highviewapps.com/exports/export-profile/"
)

$newTimestamp = 44233.00265578369
$startRow = 884

for ($i = 0; $i -lt 14; $i++) {
    $row = $startRow + $i
    $url = $urls[$i]

    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 3).Value = "Disponible"

    $dCell = $ws.Cells.Item($row, 4)
    $dCell.Value = $newTimestamp
    $dCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"

    # URL cell: text mirrors the hyperlink address (as in every prior row),
    # split into the address + optional in-page fragment ("location").
    $bCell = $ws.Cells.Item($row, 2)
    $hashPos = $url.IndexOf("#")
    if ($hashPos -ge 0) {
        $address = $url.Substring(0, $hashPos)
        $subAddress = $url.Substring($hashPos + 1)
    } else {
        $address = $url
        $subAddress = ""
    }

    $bCell.Value = $url
    $ws.Hyperlinks.Add($bCell, $address, $subAddress) | Out-Null
    # Hyperlinks.Add() re-stamps its own xf on the cell; re-apply the sheet's
    # normal "Hyperlink" cell style afterwards so it matches every other
    # hyperlinked cell in the column instead of getting a fresh duplicate style.
    $bCell.Style = "Hyperlink"
}
